# Apply the edit described by the diff:
#  - Rename the odor value "30pct_oct" -> "30pct" for every data row (B2:B59)
#  - Change the date format on column E (E2:E59) from mm-dd-yy to a custom
#    yyyy-mm-dd format
#  - Adjust column E width (auto result of the new, slightly wider format)
#  - Update the sheet view: no frozen/scrolled top-left cell, new selection B2:B59

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename odor value "30pct_oct" -> "30pct" in column B (rows 2-59)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq "30pct_oct") {
        $cell.Value2 = "30pct"
    }
}

# 2. Update date number format on column E data cells to yyyy-mm-dd
$ws.Range("E2:E59").NumberFormat = "yyyy\-mm\-dd;@"

# 3. Set column E width to match the new best-fit width for the new,
#    slightly wider date format (engine rounds to whole-character widths)
$ws.Columns.Item(5).ColumnWidth = 9

# 4. Update the selection / view state
$ws.Range("B2:B59").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
